$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style from A175 (s="2": bold font, border, centered,
# custom YYYY-MM-DD HH:MM:SS number format) so every new date cell in column A
# matches the formatting used by the existing rows.
$ws.Range("A175").Copy()

# Row 176
$ws.Range("A176").PasteSpecial(-4122)
$ws.Range("A176").Value = 45606
$ws.Range("B176").Value = 711.3372771078
$ws.Range("C176").Value = 220.7121626835
$ws.Range("I176").Value = 376.118382017
$ws.Range("K176").Value = 79.555150438797
$ws.Range("N176").Value = 38.48587173312
$ws.Range("O176").Value = 1.0680105051
$ws.Range("Q176").Value = 0.0000029568
$ws.Range("U176").Value = 289.8871550873386
$ws.Range("Z176").Value = 669.82803943252

# Row 177
$ws.Range("A177").PasteSpecial(-4122)
$ws.Range("A177").Value = 45607
$ws.Range("B177").Value = 784.6038569322001
$ws.Range("C177").Value = 233.7737442965
$ws.Range("I177").Value = 397.711125803
$ws.Range("K177").Value = 81.399626165871
$ws.Range("N177").Value = 42.22270319904001
$ws.Range("O177").Value = 1.1230002861
$ws.Range("Q177").Value = 0.0000031176
$ws.Range("U177").Value = 303.8314180637375
$ws.Range("Z177").Value = 700.80409271643

# Row 178
$ws.Range("A178").PasteSpecial(-4122)
$ws.Range("A178").Value = 45608
$ws.Range("B178").Value = 778.4438910678
$ws.Range("C178").Value = 224.91325213
$ws.Range("I178").Value = 378.929377684
$ws.Range("K178").Value = 90.91323781077899
$ws.Range("N178").Value = 39.80132953888
$ws.Range("O178").Value = 1.0646055651
$ws.Range("Q178").Value = 0.0000033912
$ws.Range("U178").Value = 267.371647896089
$ws.Range("Z178").Value = 710.912067998548

# Row 179
$ws.Range("A179").PasteSpecial(-4122)
$ws.Range("A179").Value = 45609
$ws.Range("B179").Value = 799.891012656
$ws.Range("C179").Value = 220.986041266
$ws.Range("I179").Value = 385.178024103
$ws.Range("K179").Value = 85.91373781371
$ws.Range("N179").Value = 37.37995587872
$ws.Range("O179").Value = 1.0566890796
$ws.Range("Q179").Value = 0.0000049944
$ws.Range("U179").Value = 258.1607585905778
$ws.Range("Z179").Value = 648.028350805648

# Row 180
$ws.Range("A180").PasteSpecial(-4122)
$ws.Range("A180").Value = 45610
$ws.Range("B180").Value = 772.8995854602
$ws.Range("C180").Value = 212.087414107
$ws.Range("I180").Value = 375.223160467
$ws.Range("K180").Value = 85.671043639095
$ws.Range("N180").Value = 35.63377295072
$ws.Range("O180").Value = 1.059447081
$ws.Range("Q180").Value = 0.000004956
$ws.Range("U180").Value = 261.1031260076161
$ws.Range("Z180").Value = 784.6024684423759

# Row 181
$ws.Range("A181").PasteSpecial(-4122)
$ws.Range("A181").Value = 45611
$ws.Range("B181").Value = 805.7048245146001
$ws.Range("C181").Value = 214.2500148635
$ws.Range("I181").Value = 391.122295195
$ws.Range("K181").Value = 92.806252372776
$ws.Range("N181").Value = 37.27518490304
$ws.Range("O181").Value = 1.0544588439
$ws.Range("Q181").Value = 0.0000055992
$ws.Range("U181").Value = 265.4527126241075
$ws.Range("Z181").Value = 1144.01784759824

# Row 182
$ws.Range("A182").PasteSpecial(-4122)
$ws.Range("A182").Value = 45612
$ws.Range("B182").Value = 801.7648997976
$ws.Range("C182").Value = 217.2217708245
$ws.Range("I182").Value = 386.2164811010001
$ws.Range("K182").Value = 97.17474751584599
$ws.Range("N182").Value = 37.51965051296
$ws.Range("O182").Value = 1.0601110443
$ws.Range("Q182").Value = 0.000005112
$ws.Range("U182").Value = 290.1430131236028
$ws.Range("Z182").Value = 1196.793589509022

# Row 183
$ws.Range("A183").PasteSpecial(-4122)
$ws.Range("A183").Value = 45613
$ws.Range("B183").Value = 795.2955991722001
$ws.Range("C183").Value = 213.2786126
$ws.Range("I183").Value = 425.176522957
$ws.Range("K183").Value = 96.73789800153901
$ws.Range("N183").Value = 35.48243709696
$ws.Range("O183").Value = 1.053658683
$ws.Range("Q183").Value = 0.0000050904
$ws.Range("U183").Value = 269.930228258731
$ws.Range("Z183").Value = 1971.614146388058

# Row 184
$ws.Range("A184").PasteSpecial(-4122)
$ws.Range("A184").Value = 45614
$ws.Range("B184").Value = 800.6776699824001
$ws.Range("C184").Value = 222.41714353
$ws.Range("I184").Value = 429.330350949
$ws.Range("K184").Value = 98.14552421430599
$ws.Range("N184").Value = 37.39159709824
$ws.Range("O184").Value = 1.0537097571
$ws.Range("Q184").Value = 0.0000049152
$ws.Range("U184").Value = 281.0600528362237
$ws.Range("Z184").Value = 1764.889748682806

$excel.CutCopyMode = 0
